$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.635.93'
$ws.Range("E2").Value = '  -0.67%  '

$ws.Range("D3").Value = '1.895.64'
$ws.Range("E3").Value = '  +0.95%  '

$ws.Range("E4").Value = '  +0.51%  '

$ws.Range("D5").Value = '327.22'
$ws.Range("E5").Value = '  -0.03%  '

$ws.Range("E6").Value = '  +0.39%  '

$ws.Range("D7").Value = '0.4592'

$ws.Range("D8").Value = '0.3869'
$ws.Range("E8").Value = '  -1.63%  '

$ws.Range("D9").Value = '46.69'
$ws.Range("E9").Value = '  +0.56%  '

$ws.Range("D10").Value = '0.07869'
$ws.Range("E10").Value = '  -0.69%  '

$ws.Range("D11").Value = '1.001'
$ws.Range("E11").Value = '  +2.58%  '

$ws.Range("D12").Value = '21.70'
$ws.Range("E12").Value = '  -3.09%  '

$ws.Range("D13").Value = '1.901.49'
$ws.Range("E13").Value = '  +3.28%  '

$ws.Range("D14").Value = '7.095'
$ws.Range("E14").Value = '  +1.90%  '

$ws.Range("D15").Value = '5.715'
$ws.Range("E15").Value = '  -0.74%  '

$ws.Range("D16").Value = '0.06964'
$ws.Range("E16").Value = '  -0.65%  '

$ws.Range("D17").Value = '87.30'
$ws.Range("E17").Value = '  -1.52%  '

$ws.Range("E18").Value = '  +0.41%  '

$ws.Range("D19").Value = '0.00001004'
$ws.Range("E19").Value = '  -0.92%  '

$ws.Range("D20").Value = '17.18'
$ws.Range("E20").Value = '  +1.05%  '

$ws.Range("D21").Value = '1.006'
$ws.Range("E21").Value = '  +0.38%  '

$ws.Range("D22").Value = '28.677.96'
$ws.Range("E22").Value = '  -0.54%  '

$ws.Range("D23").Value = '5.316'
$ws.Range("E23").Value = '  -0.56%  '

$ws.Range("E24").Value = '  -0.92%  '

$ws.Range("D25").Value = '2.138.06'
$ws.Range("E25").Value = '  +3.34%  '

$ws.Range("D26").Value = '2.062'
$ws.Range("E26").Value = '  -2.43%  '

$ws.Range("D27").Value = '154.83'
$ws.Range("E27").Value = '  +0.64%  '

$ws.Range("D28").Value = '19.33'
$ws.Range("E28").Value = '  -0.44%  '

$ws.Range("D29").Value = '5.841'
$ws.Range("E29").Value = '  +1.38%  '

$ws.Range("D30").Value = '118.48'
$ws.Range("E30").Value = '  -1.10%  '

$ws.Range("D31").Value = '1.932'
$ws.Range("E31").Value = '  -3.97%  '

$ws.Range("D32").Value = '0.09317'

$ws.Range("D33").Value = '0.9285'
$ws.Range("E33").Value = '  -1.39%  '

$ws.Range("D34").Value = '5.304'
$ws.Range("E34").Value = '  -0.47%  '

$ws.Range("D35").Value = '1.336'
$ws.Range("E35").Value = '  -1.38%  '

$ws.Range("D36").Value = '3.270'
$ws.Range("E36").Value = '  -2.58%  '

$ws.Range("D37").Value = '0.05748'
$ws.Range("E37").Value = '  -2.22%  '

$ws.Range("D38").Value = '1.157'
$ws.Range("E38").Value = '  +0.92%  '

$ws.Range("D39").Value = '0.02073'
$ws.Range("E39").Value = '  -2.33%  '

$ws.Range("D40").Value = '7.762'
$ws.Range("E40").Value = '  -2.14%  '

$ws.Range("D41").Value = '0.5641'

$ws.Range("E42").Value = '  -0.43%  '

$ws.Range("D43").Value = '9.757'
$ws.Range("E43").Value = '  -2.35%  '

$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").Value = '2.221'
$ws.Range("E44").Value = '  +4.02%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '11.77'
$ws.Range("E45").Value = '  +0.03%  '

$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").Value = '0.07166'
$ws.Range("E46").Value = '  -0.97%  '

$ws.Range("D47").Value = '0.5335'
$ws.Range("E47").Value = '  +0.02%  '

$ws.Range("D48").Value = '1.118'
$ws.Range("E48").Value = '  -1.47%  '

$ws.Range("D49").Value = '1.834'
$ws.Range("E49").Value = '  -1.12%  '

$ws.Range("D50").Value = '112.91'
$ws.Range("E50").Value = '  -0.95%  '

$ws.Range("D51").Value = '2.461'
$ws.Range("E51").Value = '  +4.15%  '
